# Update Leave Card 12/22/2023 10:59 AM
$wb = $excel.ActiveWorkbook

$wsLeave = $wb.Worksheets.Item("LEAVE CREDITS")
$wsConv  = $wb.Worksheets.Item("CONVERTION")

# --- LEAVE CREDITS sheet: fill in PERIOD dates for rows 12-22 (month-end dates) ---
$dates = @(45077, 45107, 45138, 45169, 45199, 45230, 45260, 45291, 45322, 45351, 45382)
$row = 12
foreach ($d in $dates) {
    $wsLeave.Cells.Item($row, 1).Value = [DateTime]::FromOADate($d)
    $row++
}

# EARNED values of 1.25 for rows 12-16 (first 5 periods) in column C
for ($r = 12; $r -le 16; $r++) {
    $wsLeave.Cells.Item($r, 3).Value = 1.25
}

# --- CONVERTION sheet: add "TOTAL LEAVE BALANCE" label and formula ---
$wsConv.Range("A6").Value = "TOTAL LEAVE BALANCE"
$wsConv.Range("A7").Formula = "=SUM('LEAVE CREDITS'!E9,'LEAVE CREDITS'!I9)"

# --- Selections / active sheet / view state ---
$wsLeave.Range("C12:C16").Select()
$wsConv.Range("A6:A7").Select()
$wsConv.Activate()

$wb.Save()
